$d = $word.ActiveDocument

# The document currently ends with an empty paragraph that only holds the
# "_GoBack" bookmark (a leftover from the last edit position). The new
# version moves that bookmark up into the "Here we see..." paragraph, right
# after the newly-inserted "over 60" text, and leaves the trailing
# paragraph completely empty.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Locate the "XX" placeholder inside the "Here we see..." paragraph.
$rng = $d.Content
$rng.Find.Execute("XX")

# Mark the start of "XX" with a temporary bookmark so that, once the text is
# replaced, the replacement text ("over 60") ends up in its own run instead
# of being merged into the preceding run.
$startRng = $rng.Duplicate
$startRng.Collapse(1)  # wdCollapseStart
$d.Bookmarks.Add("TempStart", $startRng)

# Mark the end of "XX" - this is where the real "_GoBack" bookmark needs to
# live once the replacement text is in place.
$endRng = $rng.Duplicate
$endRng.Collapse(0)  # wdCollapseEnd
$d.Bookmarks.Add("_GoBack", $endRng)

# Replace "XX" with "over 60".
$rng.Text = "over 60"

# Remove the temporary bookmark now that the run split it created will
# persist.
if ($d.Bookmarks.Exists("TempStart")) {
    $d.Bookmarks.Item("TempStart").Delete()
}
